# Weekly refresh of the "Vega Monumental Concepción - Haba" sheet:
# a new market-day row is inserted above the existing row 62, shifting
# all subsequent rows down by one (old row 62..69 -> new row 63..70).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 62 (pushes 62..69 down to 63..70,
# carrying the date number-format down with them).
$ws.Rows(62).Insert()

# Populate the newly inserted row 62 with the new weekly record.
$ws.Cells.Item(62, 1).Value = 11
$ws.Cells.Item(62, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(62, 3).Value = "Bíobío"
$ws.Cells.Item(62, 4).Value = 45218
$ws.Cells.Item(62, 5).Value = 8
$ws.Cells.Item(62, 6).Value = 100112026
$ws.Cells.Item(62, 7).Value = "Haba"
$ws.Cells.Item(62, 8).Value = "Sin especificar"
$ws.Cells.Item(62, 9).Value = "Primera"
$ws.Cells.Item(62, 10).Value = 100
$ws.Cells.Item(62, 11).Value = 11000
$ws.Cells.Item(62, 12).Value = 12000
$ws.Cells.Item(62, 13).Value = 11500
$ws.Cells.Item(62, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(62, 15).Value = "Región Metropolitana"
$ws.Cells.Item(62, 16).Value = 460
$ws.Cells.Item(62, 17).Value = 25
$ws.Cells.Item(62, 18).Value = "Hortaliza"
